$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L9").Value  = 1074.68   # ASTUDILLO ESPINOZA JOSE MANUEL - PIEDRA SINTERIZADA
$wsGrupo.Range("L12").Value = 2103.17   # BUSTAMANTE ROSERO MARCO TULIO - PIEDRA SINTERIZADA
$wsGrupo.Range("M16").Value = 3382.14   # DECOGARCIA S.A.S. - PORCELANATO
$wsGrupo.Range("M35").Value = 3057.58   # ORTEGA ROMAN LUIS FERNANDO - PORCELANATO

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F9").Value  = 3461.48   # ASTUDILLO ESPINOZA JOSE MANUEL - septiembre
$wsMensual.Range("F12").Value = 2621.57   # BUSTAMANTE ROSERO MARCO TULIO - septiembre
$wsMensual.Range("F16").Value = 3382.14   # DECOGARCIA S.A.S. - septiembre
$wsMensual.Range("F35").Value = 4132.26   # ORTEGA ROMAN LUIS FERNANDO - septiembre
$wsMensual.Range("F59").Value = 88645.14  # TOTAL - septiembre

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PIEDRA SINTERIZADA (row 11)
$wsCumpl.Range("D11").Value = 16711.27
$wsCumpl.Range("E11").Value = 1120.1443984654
$wsCumpl.Range("F11").Value = 0.9371814050509756

# PORCELANATO (row 12)
$wsCumpl.Range("D12").Value = 52203.09
$wsCumpl.Range("E12").Value = 9660.630394756605
$wsCumpl.Range("F12").Value = 0.8438401322598856

# TOTAL (row 15)
$wsCumpl.Range("D15").Value = 86282.28
$wsCumpl.Range("E15").Value = 35772.55551083435
$wsCumpl.Range("F15").Value = 0.7069140656237338
